$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range("D2")
$r.NumberFormat = "@"
$r.Value = "30.011.10"
$r.Style = "Normal"
$ws.Range("E2").Value = "  +0.08%  "

$r = $ws.Range("D3")
$r.NumberFormat = "@"
$r.Value = "1.910.25"
$r.Style = "Normal"
$ws.Range("E3").Value = "  +0.42%  "

$r = $ws.Range("D4")
$r.NumberFormat = "@"
$r.Value = "0.9997"
$r.Style = "Normal"
$ws.Range("E4").Value = "  -0.03%  "

$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "0.7932"
$r.Style = "Normal"
$ws.Range("E5").Value = "  +6.51%  "

$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = "242.13"
$r.Style = "Normal"
$ws.Range("E6").Value = "  +0.34%  "

$r = $ws.Range("D7")
$r.NumberFormat = "@"
$r.Value = "1.000"
$r.Style = "Normal"
$ws.Range("E7").Value = "  +0.00%  "

$r = $ws.Range("D8")
$r.NumberFormat = "@"
$r.Value = "0.3168"
$r.Style = "Normal"
$ws.Range("E8").Value = "  +3.20%  "

$r = $ws.Range("D9")
$r.NumberFormat = "@"
$r.Value = "26.37"
$r.Style = "Normal"
$ws.Range("E9").Value = "  +3.16%  "

$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = "0.06898"
$r.Style = "Normal"
$ws.Range("E10").Value = "  +0.10%  "

$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = "0.08002"
$r.Style = "Normal"
$ws.Range("E11").Value = "  -0.17%  "

$r = $ws.Range("D12")
$r.NumberFormat = "@"
$r.Value = "1.908.82"
$r.Style = "Normal"
$ws.Range("E12").Value = "  +0.32%  "

$r = $ws.Range("D13")
$r.NumberFormat = "@"
$r.Value = "0.7438"
$r.Style = "Normal"
$ws.Range("E13").Value = "  -1.39%  "

$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = "5.196"
$r.Style = "Normal"
$ws.Range("E14").Value = "  -1.30%  "

$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = "93.18"
$r.Style = "Normal"
$ws.Range("E15").Value = "  +1.84%  "

$r = $ws.Range("D16")
$r.NumberFormat = "@"
$r.Value = "30.010.95"
$r.Style = "Normal"
$ws.Range("E16").Value = "  +0.08%  "

$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = "13.95"
$r.Style = "Normal"
$ws.Range("E17").Value = "  -0.58%  "

$r = $ws.Range("D18")
$r.NumberFormat = "@"
$r.Value = "5.874"
$r.Style = "Normal"
$ws.Range("E18").Value = "  -4.79%  "

$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = "246.00"
$r.Style = "Normal"
$ws.Range("E19").Value = "  +3.69%  "

$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = "0.000007741"
$r.Style = "Normal"
$ws.Range("E20").Value = "  -0.01%  "

$ws.Range("E21").Value = "  -0.02%  "

$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = "2.143.53"
$r.Style = "Normal"
$ws.Range("E22").Value = "  -0.40%  "

$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = "0.9994"
$r.Style = "Normal"
$ws.Range("E23").Value = "  -0.02%  "

$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = "6.843"
$r.Style = "Normal"
$ws.Range("E24").Value = "  -3.70%  "

$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = "168.15"
$r.Style = "Normal"
$ws.Range("E25").Value = "  +1.06%  "

$r = $ws.Range("D26")
$r.NumberFormat = "@"
$r.Value = "9.236"
$r.Style = "Normal"
$ws.Range("E26").Value = "  -0.75%  "

$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = "0.1397"
$r.Style = "Normal"
$ws.Range("E27").Value = "  +10.51%  "

$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = "18.91"
$r.Style = "Normal"
$ws.Range("E28").Value = "  +0.54%  "

$ws.Range("E29").Value = "  -0.94%  "

$ws.Range("E30").Value = "  +1.63%  "

$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = "1.520"
$r.Style = "Normal"
$ws.Range("E31").Value = "  -0.23%  "

$r = $ws.Range("D32")
$r.NumberFormat = "@"
$r.Value = "4.319"
$r.Style = "Normal"
$ws.Range("E32").Value = "  +0.50%  "

$r = $ws.Range("D33")
$r.NumberFormat = "@"
$r.Value = "4.089"
$r.Style = "Normal"
$ws.Range("E33").Value = "  +1.10%  "

$r = $ws.Range("D34")
$r.NumberFormat = "@"
$r.Value = "0.05546"
$r.Style = "Normal"
$ws.Range("E34").Value = "  +2.78%  "

$r = $ws.Range("D35")
$r.NumberFormat = "@"
$r.Value = "1.256"
$r.Style = "Normal"
$ws.Range("E35").Value = "  -2.15%  "

$r = $ws.Range("D36")
$r.NumberFormat = "@"
$r.Value = "0.7341"
$r.Style = "Normal"
$ws.Range("E36").Value = "  -0.52%  "

$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = "2.720"
$r.Style = "Normal"
$ws.Range("E37").Value = "  -0.16%  "

$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = "0.01924"
$r.Style = "Normal"
$ws.Range("E38").Value = "  -0.92%  "

$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = "2.790"
$r.Style = "Normal"
$ws.Range("E39").Value = "  +0.72%  "

$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = "6.148"
$r.Style = "Normal"
$ws.Range("E40").Value = "  -1.29%  "

$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = "0.4424"
$r.Style = "Normal"
$ws.Range("E41").Value = "  -0.69%  "

$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = "72.35"
$r.Style = "Normal"
$ws.Range("E42").Value = "  -0.32%  "

$ws.Range("E43").Value = "  +0.01%  "

$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = "0.8366"
$r.Style = "Normal"
$ws.Range("E44").Value = "  +0.63%  "

$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = "1.875"
$r.Style = "Normal"
$ws.Range("E45").Value = "  -3.45%  "

$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = "100.54"
$r.Style = "Normal"
$ws.Range("E46").Value = "  -1.03%  "

$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = "7.553"
$r.Style = "Normal"
$ws.Range("E47").Value = "  -1.72%  "

$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = "989.93"
$r.Style = "Normal"

$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = "2.053.55"
$r.Style = "Normal"
$ws.Range("E49").Value = "  -0.17%  "

$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = "36.29"
$r.Style = "Normal"
$ws.Range("E50").Value = "  -0.78%  "

$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = "1.481"
$r.Style = "Normal"
$ws.Range("E51").Value = "  +0.37%  "
